# Update transition-probability matrix cells on Sheet1 with refreshed
# simulation results (more games simulated, faster simulate-game logic,
# and initial optimization-logic groundwork).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.2051671732522796
$ws.Range("C2").Value = 0.5273556231003039
$ws.Range("J2").Value = 0.01671732522796352
$ws.Range("P2").Value = 0.1595744680851064
$ws.Range("S2").Value = 0.0911854103343465
$ws.Range("C3").Value = 0.01680672268907563
$ws.Range("J3").Value = 0.02801120448179272
$ws.Range("P3").Value = 0.7310924369747899
$ws.Range("S3").Value = 0.2240896358543417
$ws.Range("J4").Value = 0.07894736842105263
$ws.Range("P4").Value = 0.6052631578947368
$ws.Range("S4").Value = 0.3157894736842105
$ws.Range("B6").Value = 0.05764966740576496
$ws.Range("D6").Value = 0.0221729490022173
$ws.Range("F6").Value = 0.07538802660753881
$ws.Range("J6").Value = 0.2195121951219512
$ws.Range("O6").Value = 0.0376940133037694
$ws.Range("Q6").Value = 0.1751662971175166
$ws.Range("R6").Value = 0.05764966740576496
$ws.Range("S6").Value = 0.3547671840354767
$ws.Range("B7").Value = 0.09498680738786279
$ws.Range("D7").Value = 0.0237467018469657
$ws.Range("F7").Value = 0.03430079155672823
$ws.Range("J7").Value = 0.1503957783641161
$ws.Range("O7").Value = 0.0316622691292876
$ws.Range("Q7").Value = 0.1952506596306069
$ws.Range("R7").Value = 0.07651715039577836
$ws.Range("S7").Value = 0.3931398416886543
$ws.Range("B8").Value = 0.1219806763285024
$ws.Range("D8").Value = 0.01207729468599034
$ws.Range("E8").Value = 0.001207729468599034
$ws.Range("F8").Value = 0.06521739130434782
$ws.Range("J8").Value = 0.1364734299516908
$ws.Range("O8").Value = 0.01811594202898551
$ws.Range("Q8").Value = 0.1835748792270532
$ws.Range("R8").Value = 0.05917874396135266
$ws.Range("S8").Value = 0.4021739130434783
$ws.Range("B9").Value = 0.0972972972972973
$ws.Range("D9").Value = 0.02432432432432433
$ws.Range("E9").Value = 0.002702702702702703
$ws.Range("F9").Value = 0.06486486486486487
$ws.Range("J9").Value = 0.1432432432432433
$ws.Range("O9").Value = 0.02162162162162162
$ws.Range("Q9").Value = 0.2027027027027027
$ws.Range("R9").Value = 0.08378378378378379
$ws.Range("S9").Value = 0.3594594594594595
$ws.Range("B10").Value = 0.1156716417910448
$ws.Range("D10").Value = 0.01604477611940298
$ws.Range("E10").Value = 0.001119402985074627
$ws.Range("F10").Value = 0.07313432835820896
$ws.Range("J10").Value = 0.133955223880597
$ws.Range("O10").Value = 0.0291044776119403
$ws.Range("Q10").Value = 0.2223880597014925
$ws.Range("R10").Value = 0.06417910447761194
$ws.Range("S10").Value = 0.3444029850746269
$ws.Range("F11").Value = 0.0015220700152207
$ws.Range("G11").Value = 0.1400304414003044
$ws.Range("J11").Value = 0.106544901065449
$ws.Range("K11").Value = 0.1948249619482496
$ws.Range("L11").Value = 0.5494672754946728
$ws.Range("S11").Value = 0.0076103500761035
$ws.Range("G12").Value = 0.6955380577427821
$ws.Range("J12").Value = 0.2178477690288714
$ws.Range("K12").Value = 0.007874015748031496
$ws.Range("L12").Value = 0.03674540682414698
$ws.Range("S12").Value = 0.04199475065616798
$ws.Range("G13").Value = 0.7017543859649122
$ws.Range("J13").Value = 0.2807017543859649
$ws.Range("S13").Value = 0.01754385964912281
$ws.Range("F15").Value = 0.018
$ws.Range("H15").Value = 0.132
$ws.Range("I15").Value = 0.076
$ws.Range("J15").Value = 0.368
$ws.Range("K15").Value = 0.074
$ws.Range("M15").Value = 0.006
$ws.Range("O15").Value = 0.052
$ws.Range("S15").Value = 0.274
$ws.Range("F16").Value = 0.02015113350125945
$ws.Range("H16").Value = 0.1662468513853904
$ws.Range("I16").Value = 0.09571788413098237
$ws.Range("J16").Value = 0.4433249370277078
$ws.Range("K16").Value = 0.1183879093198992
$ws.Range("M16").Value = 0.007556675062972292
$ws.Range("O16").Value = 0.04534005037783375
$ws.Range("S16").Value = 0.1032745591939547
$ws.Range("F17").Value = 0.0134297520661157
$ws.Range("H17").Value = 0.1776859504132231
$ws.Range("I17").Value = 0.07747933884297521
$ws.Range("J17").Value = 0.4483471074380165
$ws.Range("K17").Value = 0.09400826446280992
$ws.Range("M17").Value = 0.0134297520661157
$ws.Range("N17").Value = 0.002066115702479339
$ws.Range("O17").Value = 0.07541322314049587
$ws.Range("S17").Value = 0.0981404958677686
$ws.Range("F18").Value = 0.0198019801980198
$ws.Range("H18").Value = 0.1551155115511551
$ws.Range("I18").Value = 0.08580858085808581
$ws.Range("J18").Value = 0.4554455445544555
$ws.Range("K18").Value = 0.1254125412541254
$ws.Range("M18").Value = 0.0033003300330033
$ws.Range("O18").Value = 0.0627062706270627
$ws.Range("S18").Value = 0.0924092409240924
$ws.Range("F19").Value = 0.01685855263157895
$ws.Range("H19").Value = 0.1994243421052632
$ws.Range("I19").Value = 0.07976973684210527
$ws.Range("J19").Value = 0.3712993421052632
$ws.Range("K19").Value = 0.1233552631578947
$ws.Range("M19").Value = 0.01521381578947368
$ws.Range("N19").Value = 0.0004111842105263158
$ws.Range("O19").Value = 0.0756578947368421
$ws.Range("S19").Value = 0.1180098684210526

Write-Output "Updated 111 cells with refreshed matrix values."
